# Auto-generated COM-interop script implementing the volatile-keyword.docx edit.
$d = $word.ActiveDocument

# 1) Replace the text of the first (existing) paragraph.
$p1 = $d.Paragraphs.Item(1).Range
$p1.Text = "In Java, the volatile keyword is used as a modifier for variables to indicate that the variable's value may be modified by multiple threads, and any thread reading the variable will see the most up-to-date value. It ensures that changes made to the variable by one thread are immediately visible to other threads, preventing potential visibility and ordering issues that can occur in multi-threaded environments."

# 2) Create a placeholder paragraph after the existing (2nd, empty) paragraph;
#    we will grow the document by repeatedly injecting raw OOXML into the
#    trailing placeholder paragraph (the InsertXML host here mishandles
#    fragments of >=3 paragraphs at once, so we feed it <=2 at a time).
$r = $d.Paragraphs.Item(2).Range
$r.InsertParagraphAfter()

function Insert-RawParagraphs($fragBody) {
    $placeholder = $word.ActiveDocument.Paragraphs.Item($word.ActiveDocument.Paragraphs.Count).Range
    $full = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$fragBody
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $placeholder.InsertXML($full)
}

# 3) Insert the 26 new paragraphs, two at a time.
Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>When a variable is declared as volatile, the Java Memory Model (JMM) enforces the following behaviors:</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Visibility</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>: Any write to a volatile variable is immediately visible to all other threads that subsequently read the value. Similarly, any read of a volatile variable returns the most recently written value by any thread.</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>No Reordering</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>: The JMM prevents the reordering of volatile variable reads and writes with respect to other memory operations. This means that if a variable is updated by one thread, it will be visible to other threads before any subsequent memory operations.</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>It's important to note that using volatile is not a replacement for proper synchronization mechanisms like locks or the synchronized keyword. While volatile provides visibility guarantees for individual variable accesses, it does not provide atomicity for compound operations or ensure mutually exclusive access to shared resources.</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Here's an example of using a volatile variable:</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>java</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Copy code</w:t>
      </w:r>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">public class </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>SharedData</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> {</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">    private volatile int counter = 0;</w:t>
      </w:r>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">    public void </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>increment(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>) {</w:t>
      </w:r>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">        counter++; // This is a compound operation, not atomic.</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">    }</w:t>
      </w:r>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">    public int </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>getCounter</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>) {</w:t>
      </w:r>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">        return counter;</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">    }</w:t>
      </w:r>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>}</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>In this example, the counter variable is declared as volatile. Without the volatile keyword, there would be no guarantee that changes to counter made by one thread would be immediately visible to other threads. The volatile keyword ensures that any update to the counter variable is visible to all threads, avoiding any potential inconsistency issues.</w:t>
      </w:r>
    </w:p>
"@

Insert-RawParagraphs @"
<w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
<w:p>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>It's important to use volatile judiciously and only when it's necessary to ensure visibility of changes across threads. In many cases, the use of higher-level synchronization mechanisms, such as locks or synchronized blocks/methods, may be more appropriate for ensuring both visibility and atomicity for shared data access.</w:t>
      </w:r>
    </w:p>
"@

# 4) Drop the trailing placeholder paragraph left over from step 2.
$d.Paragraphs.Item($d.Paragraphs.Count).Range.Delete()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
